$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The author removed a duplicate data row (old row 5) from the sheet,
# which shifts all subsequent rows up by one and shrinks the used range
# from A1:G35 to A1:G34.
$ws.Rows("5:5").Delete()

# Mirror the resulting selection state left behind by an Excel
# "right-click row header -> Delete" operation: the whole row at the
# deletion point (now row 5) ends up selected.
$ws.Range("A5:XFD5").Select()
